$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 63.35160166666666
$ws.Cells.Item(2, 8).Value = 190.054805
$ws.Cells.Item(2, 9).Value = 0.05799618194253496
$ws.Cells.Item(2, 10).Value = 0.05832031690831512
$ws.Cells.Item(2, 13).Value = 121.928739
$ws.Cells.Item(2, 14).Value = 365.786217
$ws.Cells.Item(2, 15).Value = 0.2282232151508951
$ws.Cells.Item(2, 16).Value = 0.2419720431319445
$ws.Cells.Item(2, 17).Value = 7724.380904846964
$ws.Cells.Item(2, 18).Value = 69519.42814362267
$ws.Cells.Item(2, 19).Value = 0.01323607510940161
$ws.Cells.Item(2, 20).Value = 0.0141118862384075
$ws.Cells.Item(3, 7).Value = 63.35160166666666
$ws.Cells.Item(3, 8).Value = 190.054805
$ws.Cells.Item(3, 9).Value = 0.05799618194253496
$ws.Cells.Item(3, 10).Value = 0.05832031690831512
$ws.Cells.Item(3, 15).Value = 0.2768624053389947
$ws.Cells.Item(3, 16).Value = 0.2935413991166814
$ws.Cells.Item(3, 17).Value = 9370.609714951881
$ws.Cells.Item(3, 18).Value = 84335.48743456694
$ws.Cells.Item(3, 19).Value = 0.0160569624330882
$ws.Cells.Item(3, 20).Value = 0.01711942742219507
$ws.Cells.Item(4, 7).Value = 63.35160166666666
$ws.Cells.Item(4, 8).Value = 190.054805
$ws.Cells.Item(4, 9).Value = 0.05799618194253496
$ws.Cells.Item(4, 10).Value = 0.05832031690831512
$ws.Cells.Item(4, 13).Value = 83.50496933333334
$ws.Cells.Item(4, 14).Value = 250.514908
$ws.Cells.Item(4, 15).Value = 0.1563025480180701
$ws.Cells.Item(4, 16).Value = 0.1657186665504434
$ws.Cells.Item(4, 17).Value = 5290.173554392549
$ws.Cells.Item(4, 18).Value = 47611.56198953294
$ws.Cells.Item(4, 19).Value = 0.009064951012937798
$ws.Cells.Item(4, 20).Value = 0.009664765150845259
$ws.Cells.Item(5, 7).Value = 63.35160166666666
$ws.Cells.Item(5, 8).Value = 190.054805
$ws.Cells.Item(5, 9).Value = 0.05799618194253496
$ws.Cells.Item(5, 10).Value = 0.05832031690831512
$ws.Cells.Item(5, 13).Value = 91.06846250000001
$ws.Cells.Item(5, 14).Value = 182.136925
$ws.Cells.Item(5, 15).Value = 0.1704597085236707
$ws.Cells.Item(5, 16).Value = 0.1204857969594293
$ws.Cells.Item(5, 17).Value = 5769.33296069577
$ws.Cells.Item(5, 18).Value = 34615.99776417463
$ws.Cells.Item(5, 19).Value = 0.009886012269410285
$ws.Cells.Item(5, 20).Value = 0.007026769861624829
$ws.Cells.Item(6, 7).Value = 63.35160166666666
$ws.Cells.Item(6, 8).Value = 190.054805
$ws.Cells.Item(6, 9).Value = 0.05799618194253496
$ws.Cells.Item(6, 10).Value = 0.05832031690831512
$ws.Cells.Item(6, 13).Value = 89.83562999999999
$ws.Cells.Item(6, 14).Value = 269.50689
$ws.Cells.Item(6, 15).Value = 0.1681521229683693
$ws.Cells.Item(6, 16).Value = 0.1782820942415013
$ws.Cells.Item(6, 17).Value = 5691.231047234049
$ws.Cells.Item(6, 18).Value = 51221.07942510644
$ws.Cells.Item(6, 19).Value = 0.009752181117697057
$ws.Cells.Item(6, 20).Value = 0.01039746823524246
$ws.Cells.Item(7, 9).Value = 0.2870938079580828
$ws.Cells.Item(7, 10).Value = 0.2886983470587841
$ws.Cells.Item(7, 13).Value = 121.928739
$ws.Cells.Item(7, 14).Value = 365.786217
$ws.Cells.Item(7, 15).Value = 0.2282232151508951
$ws.Cells.Item(7, 16).Value = 0.2419720431319445
$ws.Cells.Item(7, 17).Value = 38237.37794133635
$ws.Cells.Item(7, 18).Value = 344136.4014720271
$ws.Cells.Item(7, 19).Value = 0.06552147190210728
$ws.Cells.Item(7, 20).Value = 0.0698569288866292
$ws.Cells.Item(8, 9).Value = 0.2870938079580828
$ws.Cells.Item(8, 10).Value = 0.2886983470587841
$ws.Cells.Item(8, 15).Value = 0.2768624053389947
$ws.Cells.Item(8, 16).Value = 0.2935413991166814
$ws.Cells.Item(8, 17).Value = 46386.57125084796
$ws.Cells.Item(8, 19).Value = 0.07948548222920621
$ws.Cells.Item(8, 20).Value = 0.08474491671830876
$ws.Cells.Item(9, 9).Value = 0.2870938079580828
$ws.Cells.Item(9, 10).Value = 0.2886983470587841
$ws.Cells.Item(9, 13).Value = 83.50496933333334
$ws.Cells.Item(9, 14).Value = 250.514908
$ws.Cells.Item(9, 15).Value = 0.1563025480180701
$ws.Cells.Item(9, 16).Value = 0.1657186665504434
$ws.Cells.Item(9, 17).Value = 26187.51820584619
$ws.Cells.Item(9, 18).Value = 235687.6638526157
$ws.Cells.Item(9, 19).Value = 0.04487349370405881
$ws.Cells.Item(9, 20).Value = 0.04784270510989882
$ws.Cells.Item(10, 9).Value = 0.2870938079580828
$ws.Cells.Item(10, 10).Value = 0.2886983470587841
$ws.Cells.Item(10, 13).Value = 91.06846250000001
$ws.Cells.Item(10, 14).Value = 182.136925
$ws.Cells.Item(10, 15).Value = 0.1704597085236707
$ws.Cells.Item(10, 16).Value = 0.1204857969594293
$ws.Cells.Item(10, 17).Value = 28559.46225400492
$ws.Cells.Item(10, 18).Value = 171356.7735240295
$ws.Cells.Item(10, 19).Value = 0.04893792682348549
$ws.Cells.Item(10, 20).Value = 0.03478405042624753
$ws.Cells.Item(11, 9).Value = 0.2870938079580828
$ws.Cells.Item(11, 10).Value = 0.2886983470587841
$ws.Cells.Item(11, 13).Value = 89.83562999999999
$ws.Cells.Item(11, 14).Value = 269.50689
$ws.Cells.Item(11, 15).Value = 0.1681521229683693
$ws.Cells.Item(11, 16).Value = 0.1782820942415013
$ws.Cells.Item(11, 17).Value = 28172.84066972967
$ws.Cells.Item(11, 18).Value = 253555.566027567
$ws.Cells.Item(11, 19).Value = 0.04827543329922493
$ws.Cells.Item(11, 20).Value = 0.0514697459176998
$ws.Cells.Item(12, 7).Value = 306.6739196666667
$ws.Cells.Item(12, 8).Value = 920.021759
$ws.Cells.Item(12, 9).Value = 0.2807492782203274
$ws.Cells.Item(12, 10).Value = 0.2823183583673431
$ws.Cells.Item(12, 13).Value = 121.928739
$ws.Cells.Item(12, 14).Value = 365.786217
$ws.Cells.Item(12, 15).Value = 0.2282232151508951
$ws.Cells.Item(12, 16).Value = 0.2419720431319445
$ws.Cells.Item(12, 17).Value = 37392.36430914397
$ws.Cells.Item(12, 18).Value = 336531.2787822956
$ws.Cells.Item(12, 19).Value = 0.06407350292673628
$ws.Cells.Item(12, 20).Value = 0.06831314998780252
$ws.Cells.Item(13, 7).Value = 306.6739196666667
$ws.Cells.Item(13, 8).Value = 920.021759
$ws.Cells.Item(13, 9).Value = 0.2807492782203274
$ws.Cells.Item(13, 10).Value = 0.2823183583673431
$ws.Cells.Item(13, 15).Value = 0.2768624053389947
$ws.Cells.Item(13, 16).Value = 0.2935413991166814
$ws.Cells.Item(13, 17).Value = 45361.46735596882
$ws.Cells.Item(13, 18).Value = 408253.2062037194
$ws.Cells.Item(13, 19).Value = 0.07772892046526647
$ws.Cells.Item(13, 20).Value = 0.08287212591147457
$ws.Cells.Item(14, 7).Value = 306.6739196666667
$ws.Cells.Item(14, 8).Value = 920.021759
$ws.Cells.Item(14, 9).Value = 0.2807492782203274
$ws.Cells.Item(14, 10).Value = 0.2823183583673431
$ws.Cells.Item(14, 13).Value = 83.50496933333334
$ws.Cells.Item(14, 14).Value = 250.514908
$ws.Cells.Item(14, 15).Value = 0.1563025480180701
$ws.Cells.Item(14, 16).Value = 0.1657186665504434
$ws.Cells.Item(14, 17).Value = 25608.79625709813
$ws.Cells.Item(14, 18).Value = 230479.1663138831
$ws.Cells.Item(14, 19).Value = 0.04388182754007123
$ws.Cells.Item(14, 20).Value = 0.04678542189134632
$ws.Cells.Item(15, 7).Value = 306.6739196666667
$ws.Cells.Item(15, 8).Value = 920.021759
$ws.Cells.Item(15, 9).Value = 0.2807492782203274
$ws.Cells.Item(15, 10).Value = 0.2823183583673431
$ws.Cells.Item(15, 13).Value = 91.06846250000001
$ws.Cells.Item(15, 14).Value = 182.136925
$ws.Cells.Item(15, 15).Value = 0.1704597085236707
$ws.Cells.Item(15, 16).Value = 0.1204857969594293
$ws.Cells.Item(15, 17).Value = 27928.32235289185
$ws.Cells.Item(15, 18).Value = 167569.9341173511
$ws.Cells.Item(15, 19).Value = 0.04785644013366794
$ws.Cells.Item(15, 20).Value = 0.03401535240416712
$ws.Cells.Item(16, 7).Value = 306.6739196666667
$ws.Cells.Item(16, 8).Value = 920.021759
$ws.Cells.Item(16, 9).Value = 0.2807492782203274
$ws.Cells.Item(16, 10).Value = 0.2823183583673431
$ws.Cells.Item(16, 13).Value = 89.83562999999999
$ws.Cells.Item(16, 14).Value = 269.50689
$ws.Cells.Item(16, 15).Value = 0.1681521229683693
$ws.Cells.Item(16, 16).Value = 0.1782820942415013
$ws.Cells.Item(16, 17).Value = 27550.24477782439
$ws.Cells.Item(16, 18).Value = 247952.2030004195
$ws.Cells.Item(16, 19).Value = 0.04720858715458541
$ws.Cells.Item(16, 20).Value = 0.0503323081725526
$ws.Cells.Item(17, 7).Value = 18.2131665
$ws.Cells.Item(17, 8).Value = 36.426333
$ws.Cells.Item(17, 9).Value = 0.01667351874766359
$ws.Cells.Item(17, 10).Value = 0.01117780360442777
$ws.Cells.Item(17, 13).Value = 121.928739
$ws.Cells.Item(17, 14).Value = 365.786217
$ws.Cells.Item(17, 15).Value = 0.2282232151508951
$ws.Cells.Item(17, 16).Value = 0.2419720431319445
$ws.Cells.Item(17, 17).Value = 2220.708424542043
$ws.Cells.Item(17, 18).Value = 13324.25054725226
$ws.Cells.Item(17, 19).Value = 0.003805284056470511
$ws.Cells.Item(17, 20).Value = 0.002704715975891002
$ws.Cells.Item(18, 7).Value = 18.2131665
$ws.Cells.Item(18, 8).Value = 36.426333
$ws.Cells.Item(18, 9).Value = 0.01667351874766359
$ws.Cells.Item(18, 10).Value = 0.01117780360442777
$ws.Cells.Item(18, 15).Value = 0.2768624053389947
$ws.Cells.Item(18, 16).Value = 0.2935413991166814
$ws.Cells.Item(18, 17).Value = 2693.988320025945
$ws.Cells.Item(18, 18).Value = 16163.92992015567
$ws.Cells.Item(18, 19).Value = 0.004616270505942965
$ws.Cells.Item(18, 20).Value = 0.003281148109095213
$ws.Cells.Item(19, 7).Value = 18.2131665
$ws.Cells.Item(19, 8).Value = 36.426333
$ws.Cells.Item(19, 9).Value = 0.01667351874766359
$ws.Cells.Item(19, 10).Value = 0.01117780360442777
$ws.Cells.Item(19, 13).Value = 83.50496933333334
$ws.Cells.Item(19, 14).Value = 250.514908
$ws.Cells.Item(19, 15).Value = 0.1563025480180701
$ws.Cells.Item(19, 16).Value = 0.1657186665504434
$ws.Cells.Item(19, 17).Value = 1520.889910045394
$ws.Cells.Item(19, 18).Value = 9125.339460272364
$ws.Cells.Item(19, 19).Value = 0.00260611346468688
$ws.Cells.Item(19, 20).Value = 0.00185237070828851
$ws.Cells.Item(20, 7).Value = 18.2131665
$ws.Cells.Item(20, 8).Value = 36.426333
$ws.Cells.Item(20, 9).Value = 0.01667351874766359
$ws.Cells.Item(20, 10).Value = 0.01117780360442777
$ws.Cells.Item(20, 13).Value = 91.06846250000001
$ws.Cells.Item(20, 14).Value = 182.136925
$ws.Cells.Item(20, 15).Value = 0.1704597085236707
$ws.Cells.Item(20, 16).Value = 0.1204857969594293
$ws.Cells.Item(20, 17).Value = 1658.645070411506
$ws.Cells.Item(20, 18).Value = 6634.580281646025
$ws.Cells.Item(20, 19).Value = 0.002842163145790695
$ws.Cells.Item(20, 20).Value = 0.001346766575535462
$ws.Cells.Item(21, 7).Value = 18.2131665
$ws.Cells.Item(21, 8).Value = 36.426333
$ws.Cells.Item(21, 9).Value = 0.01667351874766359
$ws.Cells.Item(21, 10).Value = 0.01117780360442777
$ws.Cells.Item(21, 13).Value = 89.83562999999999
$ws.Cells.Item(21, 14).Value = 269.50689
$ws.Cells.Item(21, 15).Value = 0.1681521229683693
$ws.Cells.Item(21, 16).Value = 0.1782820942415013
$ws.Cells.Item(21, 17).Value = 1636.191286822395
$ws.Cells.Item(21, 18).Value = 9817.14772093437
$ws.Cells.Item(21, 19).Value = 0.002803687574772539
$ws.Cells.Item(21, 20).Value = 0.001992802235617585
$ws.Cells.Item(22, 7).Value = 390.4979046666667
$ws.Cells.Item(22, 8).Value = 1171.493714
$ws.Cells.Item(22, 9).Value = 0.3574872131313914
$ws.Cells.Item(22, 10).Value = 0.3594851740611298
$ws.Cells.Item(22, 13).Value = 121.928739
$ws.Cells.Item(22, 14).Value = 365.786217
$ws.Cells.Item(22, 15).Value = 0.2282232151508951
$ws.Cells.Item(22, 16).Value = 0.2419720431319445
$ws.Cells.Item(22, 17).Value = 47612.91709814889
$ws.Cells.Item(22, 18).Value = 428516.25388334
$ws.Cells.Item(22, 19).Value = 0.08158688115617943
$ws.Cells.Item(22, 20).Value = 0.0869853620432143
$ws.Cells.Item(23, 7).Value = 390.4979046666667
$ws.Cells.Item(23, 8).Value = 1171.493714
$ws.Cells.Item(23, 9).Value = 0.3574872131313914
$ws.Cells.Item(23, 10).Value = 0.3594851740611298
$ws.Cells.Item(23, 15).Value = 0.2768624053389947
$ws.Cells.Item(23, 16).Value = 0.2935413991166814
$ws.Cells.Item(23, 17).Value = 57760.23593517387
$ws.Cells.Item(23, 18).Value = 519842.1234165649
$ws.Cells.Item(23, 19).Value = 0.09897476970549088
$ws.Cells.Item(23, 20).Value = 0.1055237809556078
$ws.Cells.Item(24, 7).Value = 390.4979046666667
$ws.Cells.Item(24, 8).Value = 1171.493714
$ws.Cells.Item(24, 9).Value = 0.3574872131313914
$ws.Cells.Item(24, 10).Value = 0.3594851740611298
$ws.Cells.Item(24, 13).Value = 83.50496933333334
$ws.Cells.Item(24, 14).Value = 250.514908
$ws.Cells.Item(24, 15).Value = 0.1563025480180701
$ws.Cells.Item(24, 16).Value = 0.1657186665504434
$ws.Cells.Item(24, 17).Value = 32608.51555392093
$ws.Cells.Item(24, 18).Value = 293476.6399852883
$ws.Cells.Item(24, 19).Value = 0.05587616229631535
$ws.Cells.Item(24, 20).Value = 0.05957340369006447
$ws.Cells.Item(25, 7).Value = 390.4979046666667
$ws.Cells.Item(25, 8).Value = 1171.493714
$ws.Cells.Item(25, 9).Value = 0.3574872131313914
$ws.Cells.Item(25, 10).Value = 0.3594851740611298
$ws.Cells.Item(25, 13).Value = 91.06846250000001
$ws.Cells.Item(25, 14).Value = 182.136925
$ws.Cells.Item(25, 15).Value = 0.1704597085236707
$ws.Cells.Item(25, 16).Value = 0.1204857969594293
$ws.Cells.Item(25, 17).Value = 35562.04378746491
$ws.Cells.Item(25, 18).Value = 213372.2627247895
$ws.Cells.Item(25, 19).Value = 0.06093716615131634
$ws.Cells.Item(25, 20).Value = 0.0433128576918544
$ws.Cells.Item(26, 7).Value = 390.4979046666667
$ws.Cells.Item(26, 8).Value = 1171.493714
$ws.Cells.Item(26, 9).Value = 0.3574872131313914
$ws.Cells.Item(26, 10).Value = 0.3594851740611298
$ws.Cells.Item(26, 13).Value = 89.83562999999999
$ws.Cells.Item(26, 14).Value = 269.50689
$ws.Cells.Item(26, 15).Value = 0.1681521229683693
$ws.Cells.Item(26, 16).Value = 0.1782820942415013
$ws.Cells.Item(26, 17).Value = 35080.62527940994
$ws.Cells.Item(26, 18).Value = 315725.6275146895
$ws.Cells.Item(26, 19).Value = 0.06011223382208937
$ws.Cells.Item(26, 20).Value = 0.06408976968038885
